$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ChildLocations")

# Mark the (new) last step of the testT4275 case as final.
$ws.Range("Q7").Value = "Final"
$ws.Range("R7").Value = "Final Transition Plan"

# Remove the now-superseded last step (row 8), shifting subsequent rows up.
$ws.Rows("8").Delete()

# Reflect the user's final selection/view after the edit.
$null = $ws.Range("F14").Select()
